$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.1486023186056301
$ws.Range("D2").Value = 0.8832211487428849

$ws.Range("C3").Value = 0.0964237805996284
$ws.Range("D3").Value = 0.9240570044770351

$ws.Range("C4").Value = -0.5814121964334773
$ws.Range("D4").Value = 0.5668726683198455

$ws.Range("C5").Value = -0.2404969981338094
$ws.Range("D5").Value = 0.8121719448551756

$ws.Range("C6").Value = 0.2361137323734779
$ws.Range("D6").Value = 0.8155287680290502

$ws.Range("C7").Value = -0.8172667813547102
$ws.Range("D7").Value = 0.4225388426556065

$ws.Range("C8").Value = -0.1429142270655231
$ws.Range("D8").Value = 0.887658786048712

$ws.Range("C9").Value = -0.878185748287354
$ws.Range("D9").Value = 0.3893334231400476

$ws.Range("C10").Value = -0.2867167853212338
$ws.Range("D10").Value = 0.7770120839040042

$ws.Range("C11").Value = 0.2709813970934803
$ws.Range("D11").Value = 0.7889316036150089

$wb.Save()
